$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook tracks localization handoff status for source files.
# This edit:
#   1) Updates the in-flight file's generated GUID-named artifact names and
#      the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#      timestamps (731d9710-... -> 65d22634-..., and new .xlf hashes/dates).
#   2) Adds a brand-new row for a freshly generated report/handoff file
#      (ffff0d17b411-9129-40d1-8123-b4510d2e7b3c.md) to all three sheets:
#      "Overview", "zh-cn", "de-de".
# ---------------------------------------------------------------------------

$oldGuidFile   = "731d9710-f734-40b9-a4c7-1ba3d02abd22.md"
$newGuidFile   = "65d22634-c923-46e5-98f1-4eb65b7cd2a0.md"
$newRowFile    = "ffff0d17b411-9129-40d1-8123-b4510d2e7b3c.md"

$newGuidPath   = "e2e\" + $newGuidFile
$newRowPath    = "e2e\" + $newRowFile

$zhXlf = "65d22634-c923-46e5-98f1-4eb65b7cd2a0.fa507b88808b6c0e0367ccdc38630035b784eacc.zh-cn.xlf"
$deXlf = "65d22634-c923-46e5-98f1-4eb65b7cd2a0.fa507b88808b6c0e0367ccdc38630035b784eacc.de-de.xlf"

$zhHandoffDate = "2016-08-30 09:07:24"
$deHandoffDate = "2016-08-30 09:07:37"
$overviewDate  = "2016-08-30 09:07:37"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06ffe488558d7f4b19cd2da9ce792bf1f2784bcf/e2e/"
$oldRowUrl = $baseUrl + $oldGuidFile
$newRowUrl = $baseUrl + $newRowFile

# ===========================================================================
# Sheet 1: Overview   (columns: A File Name, B Path And Name, C Extension,
#                       D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date)
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Drop all existing hyperlinks on the sheet so we can rebuild them cleanly
# (this keeps the existing rId2 relationship target untouched in spirit,
# since we recreate it with the very same Address it already had).
$ws1.Cells.Hyperlinks.Delete()

# -- update row 2 (existing file) --
$ws1.Range("A2").Value = $newGuidFile
$ws1.Range("B2").Value = $newGuidPath
$ws1.Range("G2").Value = $overviewDate

# -- add row 3 (new file) --
$ws1.Range("A3").Value = $newRowFile
$ws1.Range("B3").Value = $newRowPath
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = $overviewDate
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- re-create hyperlinks --
$ws1.Hyperlinks.Add($ws1.Range("B2"), $oldRowUrl, "", "", $newGuidPath)
$ws1.Range("B2").Font.Name = "Calibri"
$ws1.Range("B2").Font.Size = 11
$ws1.Range("B2").Font.Underline = 2
$ws1.Range("B2").Font.Color = 15570276

$ws1.Hyperlinks.Add($ws1.Range("B3"), $newRowUrl, "", "", $newRowPath)
$ws1.Range("B3").Font.Name = "Calibri"
$ws1.Range("B3").Font.Size = 11
$ws1.Range("B3").Font.Underline = 2
$ws1.Range("B3").Font.Color = 15570276

# -- resize table / dimension --
$ws1.ListObjects.Item("Overview").Resize($ws1.Range("A1:G3"))

# ===========================================================================
# Sheet 2: zh-cn   (columns: A Source File Name, B File Extension, C Status,
#   D Source Path, E Priority, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, I Latest Target File, J Latest Handback File,
#   K Latest Handback DateTime, L Reference Tokens, M To be localized,
#   N Dependency From, O Has metadata, P Error Detail)
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Cells.Hyperlinks.Delete()

# -- update row 2 (existing file) --
$ws2.Range("A2").Value = $newGuidFile
$ws2.Range("G2").Value = $zhXlf
$ws2.Range("H2").Value = $zhHandoffDate
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- add row 3 (new file) --
$ws2.Range("A3").Value = $newRowFile
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = $zhHandoffDate
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

# -- re-create hyperlinks --
$ws2.Hyperlinks.Add($ws2.Range("A2"), $oldRowUrl, "", "", $newGuidFile)
$ws2.Range("A2").Font.Name = "Calibri"
$ws2.Range("A2").Font.Size = 11
$ws2.Range("A2").Font.Underline = 2
$ws2.Range("A2").Font.Color = 15570276

$ws2.Hyperlinks.Add($ws2.Range("A3"), $newRowUrl, "", "", $newRowFile)
$ws2.Range("A3").Font.Name = "Calibri"
$ws2.Range("A3").Font.Size = 11
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = 15570276

# -- resize table / dimension --
$ws2.ListObjects.Item("zh-cn").Resize($ws2.Range("A1:P3"))

# ===========================================================================
# Sheet 3: de-de   (same column layout as zh-cn)
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Cells.Hyperlinks.Delete()

# -- update row 2 (existing file) --
$ws3.Range("A2").Value = $newGuidFile
$ws3.Range("G2").Value = $deXlf

# H2 (Latest Handoff Datetime) keeps referencing the same value as the
# Overview's "Latest HO Xliff Generate Date" (2016-08-30 09:07:37)
$ws3.Range("H2").Value = $overviewDate
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- add row 3 (new file) --
$ws3.Range("A3").Value = $newRowFile
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = $overviewDate
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

# -- re-create hyperlinks --
$ws3.Hyperlinks.Add($ws3.Range("A2"), $oldRowUrl, "", "", $newGuidFile)
$ws3.Range("A2").Font.Name = "Calibri"
$ws3.Range("A2").Font.Size = 11
$ws3.Range("A2").Font.Underline = 2
$ws3.Range("A2").Font.Color = 15570276

$ws3.Hyperlinks.Add($ws3.Range("A3"), $newRowUrl, "", "", $newRowFile)
$ws3.Range("A3").Font.Name = "Calibri"
$ws3.Range("A3").Font.Size = 11
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = 15570276

# -- resize table / dimension --
$ws3.ListObjects.Item("de-de").Resize($ws3.Range("A1:P3"))
